# chore(runtime): publish files + archive (2025-12-04 19:17:09)
# Refresh KHL probabilities tour sheet with the latest scrape: games have
# moved to 2025-12-04, matchups reshuffled, and all derived probability
# columns were recomputed for the new fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

$ws.Range("B2").Value = '2025-12-04T12:15:00'
$ws.Range("D2").Value = 'Локомотив'
$ws.Range("E2").Value = 897848
$ws.Range("F2").Value = 'https://text.khl.ru/text/897848.html'
$ws.Range("G2").Value = 1.09375
$ws.Range("H2").Value = 1.323529
$ws.Range("I2").Value = 1.535021
$ws.Range("J2").Value = 1.234139
$ws.Range("K2").Value = 1.163944
$ws.Range("L2").Value = 1.429275
$ws.Range("M2").Value = 2.417279
$ws.Range("N2").Value = 21.410046
$ws.Range("O2").Value = 27.038132
$ws.Range("P2").Value = 48.448178
$ws.Range("R2").Value = -0.2
$ws.Range("S2").Value = 0.30718
$ws.Range("T2").Value = 0.261547
$ws.Range("U2").Value = 0.431273
$ws.Range("V2").Value = 0.737476
$ws.Range("W2").Value = 0.262523
$ws.Range("X2").Value = 0.878381
$ws.Range("Y2").Value = 0.121619
$ws.Range("Z2").Value = 0.95146
$ws.Range("AA2").Value = 0.04854
$ws.Range("AB2").Value = 0.9830449999999999
$ws.Range("AC2").Value = 0.016955
$ws.Range("AD2").Value = 0.994746
$ws.Range("AE2").Value = 0.005254
$ws.Range("AF2").Value = 0.324304
$ws.Range("AG2").Value = 0.675696
$ws.Range("AH2").Value = 0.112789
$ws.Range("AI2").Value = 0.887211
$ws.Range("AJ2").Value = 0.418231
$ws.Range("AK2").Value = 0.581769
$ws.Range("AL2").Value = 0.173621
$ws.Range("AM2").Value = 0.826379
$ws.Range("AN2").Value = 0.792876
$ws.Range("AO2").Value = 0.875358
$ws.Range("B3").Value = '2025-12-04T12:30:00'
$ws.Range("D3").Value = 'СКА'
$ws.Range("E3").Value = 897849
$ws.Range("F3").Value = 'https://text.khl.ru/text/897849.html'
$ws.Range("G3").Value = 2.285128
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 2.927092
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 3.192564
$ws.Range("L3").Value = 3.463546
$ws.Range("M3").Value = 6.285128
$ws.Range("N3").Value = 30.885829
$ws.Range("O3").Value = 34.182529
$ws.Range("P3").Value = 65.068358
$ws.Range("Q3").Value = -0.06324399999999999
$ws.Range("R3").Value = 0.2
$ws.Range("S3").Value = 0.37997
$ws.Range("T3").Value = 0.157033
$ws.Range("U3").Value = 0.461572
$ws.Range("V3").Value = 0.101549
$ws.Range("W3").Value = 0.897025
$ws.Range("X3").Value = 0.206735
$ws.Range("Y3").Value = 0.791839
$ws.Range("Z3").Value = 0.346762
$ws.Range("AA3").Value = 0.651813
$ws.Range("AB3").Value = 0.5021
$ws.Range("AC3").Value = 0.496474
$ws.Range("AD3").Value = 0.649807
$ws.Range("AE3").Value = 0.348767
$ws.Range("AF3").Value = 0.827826
$ws.Range("AG3").Value = 0.172174
$ws.Range("AH3").Value = 0.618542
$ws.Range("AI3").Value = 0.381458
$ws.Range("AJ3").Value = 0.860208
$ws.Range("AK3").Value = 0.139792
$ws.Range("AL3").Value = 0.672358
$ws.Range("AM3").Value = 0.327642
$ws.Range("AN3").Value = 0.687709
$ws.Range("AO3").Value = 0.75752
$ws.Range("B4").Value = '2025-12-04T16:30:00'
$ws.Range("C4").Value = 'Авангард'
$ws.Range("D4").Value = 'ХК Сочи'
$ws.Range("E4").Value = 897847
$ws.Range("F4").Value = 'https://text.khl.ru/text/897847.html'
$ws.Range("G4").Value = 4.5
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1.404568
$ws.Range("J4").Value = 5.172122
$ws.Range("K4").Value = 4.836061
$ws.Range("L4").Value = 1.202284
$ws.Range("M4").Value = 5.5
$ws.Range("N4").Value = 40.265219
$ws.Range("O4").Value = 23.307977
$ws.Range("P4").Value = 63.573196
$ws.Range("Q4").Value = 0.2
$ws.Range("R4").Value = -0.2
$ws.Range("S4").Value = 0.896996
$ws.Range("T4").Value = 0.055477
$ws.Range("U4").Value = 0.036569
$ws.Range("V4").Value = 0.147815
$ws.Range("W4").Value = 0.8412269999999999
$ws.Range("X4").Value = 0.279957
$ws.Range("Y4").Value = 0.709085
$ws.Range("Z4").Value = 0.43954
$ws.Range("AA4").Value = 0.549501
$ws.Range("AB4").Value = 0.600144
$ws.Range("AC4").Value = 0.388898
$ws.Range("AD4").Value = 0.738684
$ws.Range("AE4").Value = 0.250358
$ws.Range("AF4").Value = 0.953672
$ws.Range("AG4").Value = 0.046328
$ws.Range("AH4").Value = 0.8608440000000001
$ws.Range("AI4").Value = 0.139156
$ws.Range("AJ4").Value = 0.338198
$ws.Range("AK4").Value = 0.661802
$ws.Range("AL4").Value = 0.121009
$ws.Range("AM4").Value = 0.878991
$ws.Range("AN4").Value = 0.977066
$ws.Range("AO4").Value = 0.19097

